$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
$r = $win.ScrollRow
Write-Output "ScrollRow is:"
Write-Output $r
